$wb = $excel.ActiveWorkbook
$wsHLR = $wb.Worksheets.Item("HLR")
$wsTC  = $wb.Worksheets.Item("Test Case")

# ---------------------------------------------------------------
# Test Case sheet: column widths
# ---------------------------------------------------------------
$wsTC.Columns.Item(5).ColumnWidth = 25.166666666666668   # -> stored width 26
$wsTC.Columns.Item(7).ColumnWidth = 11.666666666666666   # -> stored width ~12.43 (closest achievable)
$wsTC.Columns.Item(8).ColumnWidth = 13.0                 # -> stored width ~13.83 (closest achievable)

# ---------------------------------------------------------------
# Test Case sheet: row 2 data
# ---------------------------------------------------------------
$wsTC.Range("C2").Value = "check website url"
$wsTC.Range("D2").Value = "intenet and browser must be in the working condition"
$wsTC.Range("E2").Value = "1) web.whatsapp.com 2) press enter key"
$wsTC.Range("F2").Value = "NA"
$wsTC.Range("G2").Value = "to open the website properly"
$wsTC.Range("H2").Value = "as per expected result"
$wsTC.Range("I2").Value = "pass"
$wsTC.Range("E2").HorizontalAlignment = -4131
$wsTC.Rows.Item(2).RowHeight = 93.75

# ---------------------------------------------------------------
# Test Case sheet: row 3 data
# ---------------------------------------------------------------
$wsTC.Range("A3").Value = 2
$wsTC.Range("B3").Value = 2
$wsTC.Range("C3").Value = "check Qr Code"
$wsTC.Range("D3").Value = "intenet and browser must be in the working condition"
$wsTC.Range("E3").Value = "1) web.whatsapp.com 2) press enter key" + [char]10 + "3) scan the QR code"
$wsTC.Range("F3").Value = "NA"
$wsTC.Range("G3").Value = "when we scan the QR Code From the device Whatsapp open properly on desktop."
$wsTC.Range("H3").Value = "as per expected result"
$wsTC.Range("I3").Value = "pass"
$wsTC.Range("E3").HorizontalAlignment = -4131
$wsTC.Range("G3").Borders.LineStyle = -4142
$wsTC.Rows.Item(3).RowHeight = 187.5

# ---------------------------------------------------------------
# Test Case sheet: row 4 data
# ---------------------------------------------------------------
$wsTC.Range("A4").Value = 3
$wsTC.Range("B4").Value = 3
$wsTC.Range("C4").Value = "check need help to get started link"
$wsTC.Range("D4").Value = "intenet and browser must be in the working condition"
$wsTC.Range("E4").Value = "1) web.whatsapp.com 2) press enter key" + [char]10 + "3) click on need help to get started link"
$wsTC.Range("F4").Value = "NA"
$wsTC.Range("G4").Value = "while click on need help to get started link it is working properly and also had open a new page."
$wsTC.Range("H4").Value = "as per expected result"
$wsTC.Range("I4").Value = "pass"
$wsTC.Range("E4").HorizontalAlignment = -4131
$wsTC.Range("G4").Borders.LineStyle = -4142
$wsTC.Range("G4").HorizontalAlignment = -4131
$wsTC.Rows.Item(4).RowHeight = 225

# ---------------------------------------------------------------
# Test Case sheet: row 5 data
# ---------------------------------------------------------------
$wsTC.Range("A5").Value = 4
$wsTC.Range("B5").Value = 4
$wsTC.Range("C5").Value = "check video play button"
$wsTC.Range("D5").Value = "intenet and browser must be in the working condition"
$wsTC.Range("E5").Value = "1) web.whatsapp.com 2) press enter key" + [char]10 + "3) click on video play button"
$wsTC.Range("F5").Value = "NA"
$wsTC.Range("G5").Value = "while click on video play button it is working properly and also had display play the video."
$wsTC.Range("H5").Value = "as per expected result"
$wsTC.Range("I5").Value = "pass"
$wsTC.Range("E5").HorizontalAlignment = -4131
$wsTC.Range("G5").Borders.LineStyle = -4142
$wsTC.Range("G5").HorizontalAlignment = -4131
$wsTC.Rows.Item(5).RowHeight = 225

# ---------------------------------------------------------------
# Sheet view / selection / active sheet changes
# ---------------------------------------------------------------
$wsTC.Activate()
$wsTC.Range("E1").Select()

$wsHLR.Activate()
$wsHLR.Range("C9").Select()
